$wb = $excel.ActiveWorkbook
$ncSheet = $wb.Worksheets.Item("NC")

# --- Add "NC1" sheet right after "NC", cloned from NC (same header/label style) ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ncSheet.Copy($null, $last)
$nc1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$nc1.Name = "NC1"

$nc1.Range("B1").Value = "In-vehicle"
$nc1.Range("C1").Value = "At-stop"
$nc1.Range("D1").Value = "Extra"
$nc1.Range("E1").Value = "Total"
$nc1.Range("A2").Value = "No control"
$nc1.Range("B2").Value = 2101.086661275402
$nc1.Range("C2").Value = 12498.70440518066
$nc1.Range("D2").Value = 141.4698672425732
$nc1.Range("E2").Value = 14741.26093369863

# --- Add "HC" sheet right after "NC1", cloned from NC (same header/label style) ---
$last2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ncSheet.Copy($null, $last2)
$hc = $wb.Worksheets.Item($wb.Worksheets.Count)
$hc.Name = "HC"

$hc.Range("B1").Value = "In-vehicle"
$hc.Range("C1").Value = "At-stop"
$hc.Range("D1").Value = "Extra"
$hc.Range("E1").Value = "Total"
$hc.Range("A2").Value = "Holding control"
$hc.Range("B2").Value = 2598.39729042071
$hc.Range("C2").Value = 12418.93511103419
$hc.Range("D2").Value = 130.1347673362609
$hc.Range("E2").Value = 15147.46716879117

# Restore original active-sheet/tab selection (unchanged in the target diff)
$wb.Worksheets.Item(1).Activate()
